{"js": "// Update the UTENTE entity of the \"Schema Logico\" to add the new\n// `reset_token` and `first_login` attributes, just before the\n// `nome_ruolo(FK)` foreign key at the end of the UTENTE(...) line.\n//\n// Before: UTENTE(email, nome, cognome, username, password, nome_ruolo(FK))\n// After:  UTENTE(email, nome, cognome, username, password, reset_token, first_login, nome_ruolo(FK))\n\nconst body = context.document.body;\n\n// Locate the exact text that immediately precedes the insertion point.\n// This phrase only occurs once in the document (in the UTENTE paragraph),\n// so it safely/uniquely identifies where to insert the new fields.\nconst searchResults = body.search(\"nome_ruolo(FK))\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find 'nome_ruolo(FK))' in the document body.\");\n}\n\nconst target = searchResults.items[0];\ntarget.insertText(\"reset_token, first_login, \", Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "# Update the UTENTE entity of the \"Schema Logico\" to add the new\n# `reset_token` and `first_login` attributes, just before the\n# `nome_ruolo(FK)` foreign key at the end of the UTENTE(...) line.\n#\n# Before: UTENTE(email, nome, cognome, username, password, nome_ruolo(FK))\n# After:  UTENTE(email, nome, cognome, username, password, reset_token, first_login, nome_ruolo(FK))\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"nome_ruolo(FK))\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif ($found) {\n    # $range now corresponds to the matched text (\"nome_ruolo(FK))\");\n    # insert the new fields immediately before it.\n    $range.InsertBefore(\"reset_token, first_login, \")\n}\n"}
